$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several Price (column D) values look numeric (e.g. "537.52") and would
# otherwise be auto-converted to a Number by Excel's type inference, but in
# the source workbook the whole Price column is stored as text. A leading
# apostrophe forces Excel to keep such values as text (exactly like typing
# '537.52 into a cell), matching the original column's data type.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.829.77"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.635.20"
$ws.Range("E3").Value = "  +1.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.42%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'537.52"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'144.01"
$ws.Range("E6").Value = "  +3.36%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "'6.55"
$ws.Range("E9").Value = "  +2.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.62%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.46%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.73%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.106.88"
$ws.Range("E13").Value = "  +2.14%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "59.759.23"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'20.95"
$ws.Range("E15").Value = "  +2.24%  "

# Row 16 - now ShibaInu (was WrappedEther)
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000134"
$ws.Range("E16").Value = "  +1.23%  "

# Row 17 - now WrappedEther (was ShibaInu)
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.604.29"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'342.69"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.53%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'10.22"
$ws.Range("E20").Value = "  +1.43%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  -0.24%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'67.55"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - Polygon
$ws.Range("E24").Value = "  +1.87%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.53%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'7.24"
$ws.Range("E27").Value = "  +2.70%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  +4.81%  "

# Row 29 - USDe
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +4.05%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -0.05%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'18.97"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33 - Monero
$ws.Range("D33").Value = "'150.77"
$ws.Range("E33").Value = "  +1.34%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  +1.82%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +2.07%  "

# Row 36 - Fetch.AI
$ws.Range("D36").Value = "'0.839"
$ws.Range("E36").Value = "  +1.96%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -1.38%  "

# Row 38 - SuiNetwork
$ws.Range("D38").Value = "'0.823"
$ws.Range("E38").Value = "  +1.49%  "

# Row 39 - Bittensor
$ws.Range("D39").Value = "'289.23"
$ws.Range("E39").Value = "  +7.89%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  +1.74%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.02%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "'0.603"
$ws.Range("E42").Value = "  +0.88%  "

# Row 43 - WhiteBITCoin
$ws.Range("E43").Value = "  -0.48%  "

# Row 44 - Stellar
$ws.Range("D44").Value = "'0.0951"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0532"
$ws.Range("E45").Value = "  +3.75%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.966.38"
$ws.Range("E46").Value = "  +0.38%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "'0.0224"
$ws.Range("E47").Value = "  +1.59%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'18.50"
$ws.Range("E48").Value = "  +1.79%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  +2.78%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'110.73"
$ws.Range("E50").Value = "  -0.76%  "

# Row 51 - ZEEBU
$ws.Range("D51").Value = "'4.73"
$ws.Range("E51").Value = "  -0.14%  "
